$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I28").Value = 114
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 114
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 371
$ws.Range("N28").ClearContents()
$ws.Range("H98").Value = 4206284
$ws.Range("I98").Value = 4133764.2
$ws.Range("J98").Value = 5004000
$ws.Range("K98").Value = 4133764.2
$ws.Range("L98").Value = 5004000
$ws.Range("M98").Value = -4132266.2
$ws.Range("N98").Value = -5006996
$ws.Range("H112").Value = 3722626.2
$ws.Range("J112").Value = 4160389.5
$ws.Range("L112").Value = 12481168.5
$ws.Range("N112").Value = -12483384.5
$ws.Range("H113").Value = 6875.654
$ws.Range("I113").Value = 8220
$ws.Range("K113").Value = 8220
$ws.Range("M113").Value = -4966
$ws.Range("H122").Value = 4206284
$ws.Range("I122").Value = 4133764.2
$ws.Range("J122").Value = 5004000
$ws.Range("K122").Value = 12401292.6
$ws.Range("L122").Value = 15012000
$ws.Range("M122").Value = -12398842.6
$ws.Range("N122").Value = -15016900
$ws.Range("H129").Value = 1121.3182
$ws.Range("I129").Value = 708.45
$ws.Range("K129").Value = 2125.35
$ws.Range("M129").Value = 2874.65

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10136.047
$ws.Range("I32").Value = 7488.3945
$ws.Range("K32").Value = 7488.3945
$ws.Range("M32").Value = -7201.3945
$ws.Range("H61").Value = 2902.761
$ws.Range("I61").Value = 1722.6757
$ws.Range("J61").Value = 7754.222
$ws.Range("K61").Value = 1722.6757
$ws.Range("L61").Value = 7754.222
$ws.Range("M61").Value = -1510.6757
$ws.Range("N61").Value = -8178.222
$ws.Range("H74").Value = 63486.91
$ws.Range("I74").Value = 73545.21000000001
$ws.Range("J74").Value = 7160.4
$ws.Range("K74").Value = 73545.21000000001
$ws.Range("L74").Value = 7160.4
$ws.Range("M74").Value = -72671.21000000001
$ws.Range("N74").Value = -8908.4
$ws.Range("H77").Value = 63486.91
$ws.Range("I77").Value = 73545.21000000001
$ws.Range("J77").Value = 7160.4
$ws.Range("K77").Value = 367726.05
$ws.Range("L77").Value = 35802
$ws.Range("M77").Value = -363358.05
$ws.Range("N77").Value = -44538
$ws.Range("H122").Value = 2043.6316
$ws.Range("I122").Value = 2043.6316
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6130.8948
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3680.8948
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 3657.3225
$ws.Range("I132").Value = 3549.6785
$ws.Range("J132").Value = 4662
$ws.Range("K132").Value = 10649.0355
$ws.Range("L132").Value = 13986
$ws.Range("M132").Value = -8119.0355
$ws.Range("N132").Value = -19046
$ws.Range("H136").Value = 2902.761
$ws.Range("I136").Value = 1722.6757
$ws.Range("J136").Value = 7754.222
$ws.Range("K136").Value = 5168.0271
$ws.Range("L136").Value = 23262.666
$ws.Range("M136").Value = -2618.0271
$ws.Range("N136").Value = -28362.666
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1716.6428
$ws.Range("I94").Value = 1503.8889
$ws.Range("J94").Value = 2099.6
$ws.Range("K94").Value = 1503.8889
$ws.Range("L94").Value = 2099.6
$ws.Range("M94").Value = -1052.8889
$ws.Range("N94").Value = -3001.6
$ws.Range("H134").Value = 4002.5264
$ws.Range("J134").Value = 4425.5713
$ws.Range("L134").Value = 13276.7139
$ws.Range("N134").Value = -18346.7139

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 2160
$ws.Range("I4").Value = 50
$ws.Range("J4").Value = 8490
$ws.Range("K4").Value = 50
$ws.Range("L4").Value = 8490
$ws.Range("M4").Value = 62
$ws.Range("N4").Value = -8714
$ws.Range("H43").Value = 14411
$ws.Range("J43").Value = 14411
$ws.Range("L43").Value = 14411
$ws.Range("N43").Value = -14779
$ws.Range("H58").Value = 2127.6
$ws.Range("I58").Value = 2176.4211
$ws.Range("K58").Value = 2176.4211
$ws.Range("M58").Value = -1973.4211
$ws.Range("H95").Value = 21811.75
$ws.Range("J95").Value = 12415.667
$ws.Range("L95").Value = 12415.667
$ws.Range("N95").Value = -17907.667
$ws.Range("H101").Value = 14411
$ws.Range("J101").Value = 14411
$ws.Range("L101").Value = 14411
$ws.Range("N101").Value = -20901
$ws.Range("H132").Value = 1611.36
$ws.Range("I132").Value = 1608.0435
$ws.Range("K132").Value = 4824.1305
$ws.Range("M132").Value = -2294.1305
$ws.Range("H136").Value = 2127.6
$ws.Range("I136").Value = 2176.4211
$ws.Range("K136").Value = 6529.263300000001
$ws.Range("M136").Value = -3979.263300000001

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 1580
$ws.Range("I57").Value = 1940
$ws.Range("J57").Value = 500
$ws.Range("K57").Value = 5820
$ws.Range("L57").Value = 1500
$ws.Range("M57").Value = -5261
$ws.Range("N57").Value = -2618
$ws.Range("H107").Value = 431.22223
$ws.Range("J107").Value = 364.69232
$ws.Range("L107").Value = 1094.07696
$ws.Range("N107").Value = -4934.07696

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10013
$ws.Range("I70").Value = 10641.25
$ws.Range("J70").Value = 7500
$ws.Range("K70").Value = 10641.25
$ws.Range("L70").Value = 7500
$ws.Range("M70").Value = -10371.25
$ws.Range("N70").Value = -8040
$ws.Range("H73").Value = 10013
$ws.Range("I73").Value = 10641.25
$ws.Range("J73").Value = 7500
$ws.Range("K73").Value = 10641.25
$ws.Range("L73").Value = 7500
$ws.Range("M73").Value = -9705.25
$ws.Range("N73").Value = -9372
$ws.Range("H114").Value = 138999.5
$ws.Range("J114").Value = 138999.5
$ws.Range("L114").Value = 138999.5
$ws.Range("N114").Value = -147677.5
$ws.Range("H132").Value = 23836.969
$ws.Range("I132").Value = 28797.94
$ws.Range("J132").Value = 7631.1333
$ws.Range("K132").Value = 86393.81999999999
$ws.Range("L132").Value = 22893.3999
$ws.Range("M132").Value = -83863.81999999999
$ws.Range("N132").Value = -27953.3999

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 1461.5385
$ws.Range("I2").Value = 1461.5385
$ws.Range("K2").Value = 1461.5385
$ws.Range("M2").Value = -1349.5385
$ws.Range("H132").Value = 4104.403
$ws.Range("I132").Value = 2651.111
$ws.Range("J132").Value = 6526.5557
$ws.Range("K132").Value = 7953.333
$ws.Range("L132").Value = 19579.6671
$ws.Range("M132").Value = -5423.333
$ws.Range("N132").Value = -24639.6671

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 9998
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 9998
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 9998
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -10584
$ws.Range("H132").Value = 1551.8628
$ws.Range("I132").Value = 1449.0667
$ws.Range("K132").Value = 4347.2001
$ws.Range("M132").Value = -1817.2001
